$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being treated as text (matches the
# original inline-string cell type) instead of being auto-converted to a
# number by Excel when we assign numeric-looking values like "239.15".
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "40.938.18"
$ws.Range("E2").Value = "  -6.62%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.187.84"
$ws.Range("E3").Value = "  -7.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "239.15"
$ws.Range("E5").Value = "  -0.38%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -7.55%  "

# Row 7 - Solana
$ws.Range("D7").Value = "69.88"
$ws.Range("E7").Value = "  -5.01%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.17%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -10.52%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "36.48"
$ws.Range("E10").Value = "  +5.16%  "

# Row 11 - was Dogecoin, becomes OKB
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "57.61"
$ws.Range("E11").Value = "  -5.39%  "

# Row 12 - was OKB, becomes Dogecoin
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0937"
$ws.Range("E12").Value = "  -8.79%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -4.47%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.52"
$ws.Range("E14").Value = "  -9.81%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.515.13"
$ws.Range("E15").Value = "  -7.27%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  -10.27%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.829"
$ws.Range("E17").Value = "  -8.91%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.187.89"
$ws.Range("E18").Value = "  -7.21%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "40.903.93"
$ws.Range("E19").Value = "  -6.64%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  -9.47%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "72.39"
$ws.Range("E21").Value = "  -6.75%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  -7.83%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "230.15"
$ws.Range("E23").Value = "  -9.08%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +6.39%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.06%  "

# Row 26 - WEMIXToken
$ws.Range("E26").Value = "  -5.05%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "2.38"
$ws.Range("E27").Value = "  -4.44%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -5.15%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  -7.77%  "

# Row 30 - Monero
$ws.Range("D30").Value = "168.57"
$ws.Range("E30").Value = "  -4.19%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "20.11"
$ws.Range("E31").Value = "  -9.77%  "

# Row 32 - Kaspa
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  -9.51%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  -8.04%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.0697"
$ws.Range("E34").Value = "  -6.71%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("D35").Value = "5.04"
$ws.Range("E35").Value = "  -5.46%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -9.99%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +0.43%  "

# Row 38 - InjectiveProtocol
$ws.Range("D38").Value = "23.01"
$ws.Range("E38").Value = "  +14.29%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -7.11%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "0.0265"
$ws.Range("E40").Value = "  -4.17%  "

# Row 41 - THORChain
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  -12.08%  "

# Row 42 - MultiversX
$ws.Range("D42").Value = "63.95"
$ws.Range("E42").Value = "  -0.73%  "

# Row 43 - FTXToken
$ws.Range("D43").Value = "4.78"
$ws.Range("E43").Value = "  -11.41%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "8.62"
$ws.Range("E44").Value = "  -4.68%  "

# Row 45 - Algorand
$ws.Range("D45").Value = "0.191"
$ws.Range("E45").Value = "  -5.89%  "

# Row 46 - BinanceUSD
$ws.Range("E46").Value = "  +0.12%  "

# Row 47 - Cronos
$ws.Range("D47").Value = "0.0977"
$ws.Range("E47").Value = "  -7.83%  "

# Row 48 - SynthetixNetwork
$ws.Range("D48").Value = "4.46"
$ws.Range("E48").Value = "  +2.64%  "

# Row 49 - Celestia
$ws.Range("D49").Value = "10.15"
$ws.Range("E49").Value = "  +6.38%  "

# Row 50 - TrustWalletToken
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  -6.16%  "

# Row 51 - ARBITRUM
$ws.Range("E51").Value = "  -6.80%  "
